$wb = $excel.ActiveWorkbook

# --- Matriz_Resultados ---
$ws = $wb.Worksheets.Item("Matriz_Resultados")
$ws.Cells.Item(2,9).Value = [double]"0"
$ws.Cells.Item(3,5).Value = [double]"0"
$ws.Cells.Item(4,9).Value = [double]"0"
$ws.Cells.Item(5,3).Value = [double]"0"
$ws.Cells.Item(9,2).Value = [double]"0"
$ws.Cells.Item(9,4).Value = [double]"0"

# --- P_valores ---
$ws = $wb.Worksheets.Item("P_valores")
$ws.Cells.Item(2,3).Value = [double]"0.04498377928434749"
$ws.Cells.Item(2,4).Value = [double]"0.01307560961029575"
$ws.Cells.Item(2,5).Value = [double]"0.0004971790197476622"
$ws.Cells.Item(2,6).Value = [double]"0.0001075149407823162"
$ws.Cells.Item(2,7).Value = [double]"0.001092814935505171"
$ws.Cells.Item(2,8).Value = [double]"0.00469216381914217"
$ws.Cells.Item(2,9).Value = [double]"0.007235779763046235"
$ws.Cells.Item(2,10).Value = [double]"2.412457789091604E-09"
$ws.Cells.Item(3,2).Value = [double]"0.04498377928434749"
$ws.Cells.Item(3,4).Value = [double]"0.01476387208284291"
$ws.Cells.Item(3,5).Value = [double]"0.001408895489968032"
$ws.Cells.Item(3,6).Value = [double]"0.0005255341541874348"
$ws.Cells.Item(3,7).Value = [double]"1.148307962450801E-05"
$ws.Cells.Item(3,8).Value = [double]"5.007522925204455E-07"
$ws.Cells.Item(3,9).Value = [double]"0.02721801851139038"
$ws.Cells.Item(3,10).Value = [double]"1.707742613987762E-09"
$ws.Cells.Item(4,2).Value = [double]"0.01307560961029575"
$ws.Cells.Item(4,3).Value = [double]"0.01476387208284291"
$ws.Cells.Item(4,5).Value = [double]"0.04512865414841927"
$ws.Cells.Item(4,6).Value = [double]"0.1527187481202603"
$ws.Cells.Item(4,7).Value = [double]"0.1379302789669306"
$ws.Cells.Item(4,8).Value = [double]"0.5670572994123195"
$ws.Cells.Item(4,9).Value = [double]"0.005638392614844001"
$ws.Cells.Item(4,10).Value = [double]"2.065157544972607E-07"
$ws.Cells.Item(5,2).Value = [double]"0.0004971790197476622"
$ws.Cells.Item(5,3).Value = [double]"0.001408895489968032"
$ws.Cells.Item(5,4).Value = [double]"0.04512865414841927"
$ws.Cells.Item(5,6).Value = [double]"0.6422270523020899"
$ws.Cells.Item(5,7).Value = [double]"0.6109195650549064"
$ws.Cells.Item(5,8).Value = [double]"0.5411261269852985"
$ws.Cells.Item(5,9).Value = [double]"0.0004199423600950158"
$ws.Cells.Item(5,10).Value = [double]"3.779772401735215E-07"
$ws.Cells.Item(6,2).Value = [double]"0.0001075149407823162"
$ws.Cells.Item(6,3).Value = [double]"0.0005255341541874348"
$ws.Cells.Item(6,4).Value = [double]"0.1527187481202603"
$ws.Cells.Item(6,5).Value = [double]"0.6422270523020899"
$ws.Cells.Item(6,7).Value = [double]"0.4482606319706171"
$ws.Cells.Item(6,8).Value = [double]"0.635234240063447"
$ws.Cells.Item(6,9).Value = [double]"0.000327040347535057"
$ws.Cells.Item(6,10).Value = [double]"7.035414695266695E-08"
$ws.Cells.Item(7,2).Value = [double]"0.001092814935505171"
$ws.Cells.Item(7,3).Value = [double]"1.148307962450801E-05"
$ws.Cells.Item(7,4).Value = [double]"0.1379302789669306"
$ws.Cells.Item(7,5).Value = [double]"0.6109195650549064"
$ws.Cells.Item(7,6).Value = [double]"0.4482606319706171"
$ws.Cells.Item(7,8).Value = [double]"0.1366559950239505"
$ws.Cells.Item(7,9).Value = [double]"4.49205363350913E-07"
$ws.Cells.Item(7,10).Value = [double]"9.369675295634039E-07"
$ws.Cells.Item(8,2).Value = [double]"0.00469216381914217"
$ws.Cells.Item(8,3).Value = [double]"5.007522925204455E-07"
$ws.Cells.Item(8,4).Value = [double]"0.5670572994123195"
$ws.Cells.Item(8,5).Value = [double]"0.5411261269852985"
$ws.Cells.Item(8,6).Value = [double]"0.635234240063447"
$ws.Cells.Item(8,7).Value = [double]"0.1366559950239505"
$ws.Cells.Item(8,9).Value = [double]"1.271757561260856E-06"
$ws.Cells.Item(8,10).Value = [double]"1.335204709018711E-07"
$ws.Cells.Item(9,2).Value = [double]"0.007235779763046235"
$ws.Cells.Item(9,3).Value = [double]"0.02721801851139038"
$ws.Cells.Item(9,4).Value = [double]"0.005638392614844001"
$ws.Cells.Item(9,5).Value = [double]"0.0004199423600950158"
$ws.Cells.Item(9,6).Value = [double]"0.000327040347535057"
$ws.Cells.Item(9,7).Value = [double]"4.49205363350913E-07"
$ws.Cells.Item(9,8).Value = [double]"1.271757561260856E-06"
$ws.Cells.Item(9,10).Value = [double]"2.256623110596934E-09"
$ws.Cells.Item(10,2).Value = [double]"2.412457789091604E-09"
$ws.Cells.Item(10,3).Value = [double]"1.707742613987762E-09"
$ws.Cells.Item(10,4).Value = [double]"2.065157544972607E-07"
$ws.Cells.Item(10,5).Value = [double]"3.779772401735215E-07"
$ws.Cells.Item(10,6).Value = [double]"7.035414695266695E-08"
$ws.Cells.Item(10,7).Value = [double]"9.369675295634039E-07"
$ws.Cells.Item(10,8).Value = [double]"1.335204709018711E-07"
$ws.Cells.Item(10,9).Value = [double]"2.256623110596934E-09"

# --- Estadisticos_DM ---
$ws = $wb.Worksheets.Item("Estadisticos_DM")
$ws.Cells.Item(2,3).Value = [double]"2.201379963148599"
$ws.Cells.Item(2,4).Value = [double]"-2.841171320911328"
$ws.Cells.Item(2,5).Value = [double]"-4.502111337712737"
$ws.Cells.Item(2,6).Value = [double]"-5.323367558810774"
$ws.Cells.Item(2,7).Value = [double]"-4.094970504841406"
$ws.Cells.Item(2,8).Value = [double]"-3.357632941259305"
$ws.Cells.Item(2,9).Value = [double]"3.139876373036356"
$ws.Cells.Item(2,10).Value = [double]"-13.3246751834285"
$ws.Cells.Item(3,2).Value = [double]"-2.201379963148599"
$ws.Cells.Item(3,4).Value = [double]"-2.779500723094411"
$ws.Cells.Item(3,5).Value = [double]"-3.965258845578426"
$ws.Cells.Item(3,6).Value = [double]"-4.473152209334574"
$ws.Cells.Item(3,7).Value = [double]"-6.621497706459626"
$ws.Cells.Item(3,8).Value = [double]"-8.71198489908131"
$ws.Cells.Item(3,9).Value = [double]"2.465525480940137"
$ws.Cells.Item(3,10).Value = [double]"-13.68175999883003"
$ws.Cells.Item(4,2).Value = [double]"2.841171320911328"
$ws.Cells.Item(4,3).Value = [double]"2.779500723094411"
$ws.Cells.Item(4,5).Value = [double]"-2.199664642348018"
$ws.Cells.Item(4,6).Value = [double]"-1.512223822266675"
$ws.Cells.Item(4,7).Value = [double]"-1.57347447131931"
$ws.Cells.Item(4,8).Value = [double]"-0.5862237149779225"
$ws.Cells.Item(4,9).Value = [double]"3.265303663541629"
$ws.Cells.Item(4,10).Value = [double]"-9.374987385440917"
$ws.Cells.Item(5,2).Value = [double]"4.502111337712737"
$ws.Cells.Item(5,3).Value = [double]"3.965258845578426"
$ws.Cells.Item(5,4).Value = [double]"2.199664642348018"
$ws.Cells.Item(5,6).Value = [double]"0.4748345073810061"
$ws.Cells.Item(5,7).Value = [double]"-0.5203953785015708"
$ws.Cells.Item(5,8).Value = [double]"0.6264096423546874"
$ws.Cells.Item(5,9).Value = [double]"4.590559299016785"
$ws.Cells.Item(5,10).Value = [double]"-8.918714627837465"
$ws.Cells.Item(6,2).Value = [double]"5.323367558810774"
$ws.Cells.Item(6,3).Value = [double]"4.473152209334574"
$ws.Cells.Item(6,4).Value = [double]"1.512223822266675"
$ws.Cells.Item(6,5).Value = [double]"-0.4748345073810061"
$ws.Cells.Item(6,7).Value = [double]"-0.7802056086325619"
$ws.Cells.Item(6,8).Value = [double]"0.4849178087389195"
$ws.Cells.Item(6,9).Value = [double]"4.722411490352282"
$ws.Cells.Item(6,10).Value = [double]"-10.23091114193982"
$ws.Cells.Item(7,2).Value = [double]"4.094970504841406"
$ws.Cells.Item(7,3).Value = [double]"6.621497706459626"
$ws.Cells.Item(7,4).Value = [double]"1.57347447131931"
$ws.Cells.Item(7,5).Value = [double]"0.5203953785015708"
$ws.Cells.Item(7,6).Value = [double]"0.7802056086325619"
$ws.Cells.Item(7,8).Value = [double]"1.579010394860475"
$ws.Cells.Item(7,9).Value = [double]"8.791414263961295"
$ws.Cells.Item(7,10).Value = [double]"-8.263691138353639"
$ws.Cells.Item(8,2).Value = [double]"3.357632941259305"
$ws.Cells.Item(8,3).Value = [double]"8.71198489908131"
$ws.Cells.Item(8,4).Value = [double]"0.5862237149779225"
$ws.Cells.Item(8,5).Value = [double]"-0.6264096423546874"
$ws.Cells.Item(8,6).Value = [double]"-0.4849178087389195"
$ws.Cells.Item(8,7).Value = [double]"-1.579010394860475"
$ws.Cells.Item(8,9).Value = [double]"8.050973372575127"
$ws.Cells.Item(8,10).Value = [double]"-9.714787968350356"
$ws.Cells.Item(9,2).Value = [double]"-3.139876373036356"
$ws.Cells.Item(9,3).Value = [double]"-2.465525480940137"
$ws.Cells.Item(9,4).Value = [double]"-3.265303663541629"
$ws.Cells.Item(9,5).Value = [double]"-4.590559299016785"
$ws.Cells.Item(9,6).Value = [double]"-4.722411490352282"
$ws.Cells.Item(9,7).Value = [double]"-8.791414263961295"
$ws.Cells.Item(9,8).Value = [double]"-8.050973372575127"
$ws.Cells.Item(9,10).Value = [double]"-13.39305964194965"
$ws.Cells.Item(10,2).Value = [double]"13.3246751834285"
$ws.Cells.Item(10,3).Value = [double]"13.68175999883003"
$ws.Cells.Item(10,4).Value = [double]"9.374987385440917"
$ws.Cells.Item(10,5).Value = [double]"8.918714627837465"
$ws.Cells.Item(10,6).Value = [double]"10.23091114193982"
$ws.Cells.Item(10,7).Value = [double]"8.263691138353639"
$ws.Cells.Item(10,8).Value = [double]"9.714787968350356"
$ws.Cells.Item(10,9).Value = [double]"13.39305964194965"

# --- Resumen ---
$ws = $wb.Worksheets.Item("Resumen")
$ws.Cells.Item(2,2).Value = [double]"5"
$ws.Cells.Item(2,4).Value = [double]"3"
$ws.Cells.Item(2,5).Value = [double]"62.5"
$ws.Cells.Item(3,1).Value = "Block Bootstrapping"
$ws.Cells.Item(3,2).Value = [double]"4"
$ws.Cells.Item(3,4).Value = [double]"4"
$ws.Cells.Item(3,5).Value = [double]"50"
$ws.Cells.Item(3,6).Value = [double]"0.6645563381915951"
$ws.Cells.Item(4,1).Value = "Sieve Bootstrap"
$ws.Cells.Item(4,3).Value = [double]"0"
$ws.Cells.Item(4,4).Value = [double]"4"
$ws.Cells.Item(4,6).Value = [double]"0.6387781983476837"
$ws.Cells.Item(5,3).Value = [double]"2"
$ws.Cells.Item(5,4).Value = [double]"5"
$ws.Cells.Item(6,3).Value = [double]"0"
$ws.Cells.Item(6,4).Value = [double]"7"

